$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sexual behavior")
$ws.Range("A1").Value = "test"
Write-Host "done"
